$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Three "last modified" timestamps refreshed elsewhere in the table
#    (AchAuthLog / EmpDeductDtl / CdBranch rows) — column E holds the
#    "last modified time" value for each GenTable row.
$ws.Range("E65").Value  = "2022年03月14日 11:38:54"
$ws.Range("E77").Value  = "2022年03月14日 15:41:53"
$ws.Range("E152").Value = "2022年03月14日 15:53:02"

# 2) A new row for table "TxArchiveTableLog" is inserted right after the
#    existing "TxArchiveTable" row (row 329), pushing every following
#    row down by one.
$ws.Rows("330:330").Insert()

$ws.Range("A330").Value = "XX-系統"
$ws.Range("B330").Value = "TxArchiveTableLog"
$ws.Range("C330").Value = "歷史封存表紀錄檔"
$ws.Range("D330").Formula = "=HYPERLINK(""[\\192.168.10.16\St1Share(NAS)\SKL\DB\GenTables\XX-系統\TxArchiveTableLog.xlsx]DBD!A1"", ""連結"")"
$ws.Range("E330").Value = "2022年03月14日 15:46:04"
